# Refresh the crypto price/volume snapshot on Sheet1.
# Source: "Updated symbol list on Sun Jan 22 19:31:23 UTC 2023 with GitHub Actions"
#
# The sheet stores Price (col D) and Volume(1h) (col E) as plain text cells
# (numeric- and percent-looking strings), not real numbers. A leading
# apostrophe forces Excel to keep them as text instead of auto-converting
# them to Number/Percentage, matching how the source data is stored.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @{ D = new Price; E = new Volume(1h) }  (only the columns that changed)
$updates = [ordered]@{
    2  = @{ D = "305.02";     E = "-0.02%" }
    3  = @{ D = "36.66";      E = "2.63%" }
    4  = @{ D = "5.021";      E = "-1.61%" }
    5  = @{ D = "0.07842";    E = "-0.02%" }
    6  = @{ D = "2.165";      E = "-3.80%" }
    7  = @{ D = "8.038";      E = "-0.86%" }
    8  = @{ D = "0.9217";     E = "-0.48%" }
    9  = @{ D = "0.09952";    E = "1.02%" }
    10 = @{              E = "3.10%" }
    11 = @{ D = "0.08747";    E = "0.15%" }
    12 = @{ D = "0.03612";    E = "5.60%" }
    13 = @{ D = "0.09932";    E = "-0.16%" }
    14 = @{ D = "0.001491";   E = "0.39%" }
    15 = @{ D = "0.005624";   E = "-1.98%" }
    16 = @{ D = "3.466";      E = "-0.47%" }
    17 = @{              E = "1.28%" }
    18 = @{ D = "2.335";      E = "9.82%" }
    19 = @{              E = "0.53%" }
    20 = @{              E = "1.92%" }
    21 = @{ D = "4.920";      E = "8.00%" }
    22 = @{              E = "-1.50%" }
    23 = @{ D = "0.04616";    E = "-1.48%" }
    24 = @{ D = "0.005188";   E = "15.22%" }
    25 = @{              E = "-0.56%" }
    26 = @{ D = "0.0001401";  E = "7.81%" }
    27 = @{ D = "0.0002719";  E = "0.61%" }
    39 = @{ D = "0.01816";    E = "3.14%" }
    40 = @{ D = "0.04749";    E = "0.66%" }
    41 = @{ D = "0.007909";   E = "-1.89%" }
    42 = @{ D = "0.1406";     E = "-1.24%" }
    43 = @{ D = "0.007598";   E = "-10.86%" }
    44 = @{ D = "0.002181";   E = "-1.34%" }
    45 = @{ D = "0.01008";    E = "10.23%" }
    46 = @{ D = "0.00006354"; E = "2.28%" }
    47 = @{              E = "-0.05%" }
    48 = @{ D = "0.0005805";  E = "0.07%" }
    49 = @{              E = "805.34%" }
    50 = @{              E = "-0.01%" }
    51 = @{ D = "0.00002102"; E = "-0.05%" }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = "'" + $cols[$col]
    }
}
